# burgerzaken.docx edit script
# - merge the "<topic>" / ", toets N" run pairs into single runs
# - nudge a handful of table indents / cell margins
# - add a default header ("burgerzaken.docx") and grow the top margin to
#   make room for it
# - bump the cached PAGE field result in the footer (6 -> 7) to match the
#   now-taller layout

$d = $word.ActiveDocument

# --- 1. merge split runs in the "toets" bullet list -------------------
$pairs = @(
  "Coronavirus, toets 1",
  "Afspraken stad, toets 2",
  "Werk en Inkomen, toets 3",
  "Parkeren, toets 4"
)
foreach ($full in $pairs) {
  $d.Content.Find.Execute($full, $false, $false, $false, $false, $false, `
                           $true, 1, $false, $full, 2) | Out-Null
}

# --- 2. table indent / cell-margin tweaks ------------------------------
# tables 1-3: tblInd -22 -> -24 (dxa); dxa/20 = points
foreach ($i in 1, 2, 3) {
  $t = $d.Tables.Item($i)
  $t.Rows.LeftIndent = -24 / 20.0
}

# table 4: tblInd -81 -> -91
$t4 = $d.Tables.Item(4)
$t4.Rows.LeftIndent = -91 / 20.0

# table 5: tblInd 19 -> 17, left cell margin 6 -> 3
$t5 = $d.Tables.Item(5)
$t5.Rows.LeftIndent = 17 / 20.0
$t5.LeftPadding = 3 / 20.0
foreach ($r in $t5.Rows) {
  foreach ($c in $r.Cells) {
    $c.LeftPadding = 3 / 20.0
  }
}

# tables 6-7: tblInd 38 -> 36, left cell margin 33 -> 30
foreach ($i in 6, 7) {
  $t = $d.Tables.Item($i)
  $t.Rows.LeftIndent = 36 / 20.0
  $t.LeftPadding = 30 / 20.0
  foreach ($r in $t.Rows) {
    foreach ($c in $r.Cells) {
      $c.LeftPadding = 30 / 20.0
    }
  }
}

# --- 3. add a default header -------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdr.Range.InsertAfter("burgerzaken.docx")
$hdr.Range.Paragraphs.Item(1).Style = "Header"

# grow top margin / header distance so the header has room
$ps = $d.PageSetup
$ps.HeaderDistance = 1134 / 20.0
$ps.TopMargin = 1693 / 20.0

# --- 4. fix the cached page-count field in the footer -------------------
$ftr = $sec.Footers.Item(1)
$ftr.Range.Find.Execute("6", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "7", 2) | Out-Null

Write-Output "done"
